$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C7").Value = "Complete"
$ws.Range("C8").Value = "Complete"
$ws.Range("C9").Value = "Complete"
$ws.Range("C12").Value = "Complete"

$ws.Range("C13").Select()
